# Update the "Source" sheet test rows (35 & 36) with corrected sample data
# and adjust the current view/selection, mirroring the recorded execution log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source")

# Row 35: rename "Test1" / "test1@example.com" -> "DateTimeFormat" / "datetimeformat@example.com"
$ws.Range("B35").Value = "DateTimeFormat"
$ws.Range("D35").Value = "datetimeformat@example.com"

# Row 36: rename "Test2" / "test2@example.com" -> "Incorrect Causality" / "incorrectcausality@example.com"
$ws.Range("B36").Value = "Incorrect Causality"
$ws.Range("D36").Value = "incorrectcausality@example.com"

# Update the sheet view / selection state to match the saved workbook view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J31").Select()

